# Update the accession loader's header row to the new "vavilov accession
# format" column names (CODE / COLLECTING CODE instead of
# Accession / CollectingNumber).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CODE"
$ws.Range("B1").Value = "COLLECTING CODE"

# Move the active selection back to A1 (it was sitting on B1).
$ws.Range("A1").Select() | Out-Null
